{"js": "// Apply the two text edits from the commit:\n//  1. \"Manuel Dias\" -> \"Manuel In\u00e1cio Veladas Dias\" (representative's name)\n//  2. \"programas\" -> \"programa\" (fixes a duplicated \"s\" before the next run\n//     which already starts with \"s, projetos e a\u00e7\u00f5es...\")\nconst body = context.document.body;\n\nconst nameResults = body.search(\"Manuel Dias\", { matchCase: true, matchWholeWord: false });\nnameResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < nameResults.items.length; i++) {\n  nameResults.items[i].insertText(\"Manuel In\u00e1cio Veladas Dias\", \"Replace\");\n}\nawait context.sync();\n\nconst programResults = body.search(\"programas\", { matchCase: true, matchWholeWord: false });\nprogramResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < programResults.items.length; i++) {\n  programResults.items[i].insertText(\"programa\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Apply the two text edits from the commit:\n#  1. \"Manuel Dias\" -> \"Manuel In\u00e1cio Veladas Dias\" (representative's name)\n#  2. \"programas\" -> \"programa\" (fixes a duplicated \"s\" before the next run\n#     which already starts with \"s, projetos e a\u00e7\u00f5es...\")\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Manuel Dias\"\n$find1.Replacement.Text = \"Manuel In\u00e1cio Veladas Dias\"\n$find1.Forward = $true\n$find1.Wrap = 1\n$find1.Format = $false\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\n$find1.MatchWildcards = $false\n$find1.Execute($find1.Text, $find1.MatchCase, $find1.MatchWholeWord, $find1.MatchWildcards, $null, $null, $find1.Forward, $find1.Wrap, $find1.Format, $find1.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"programas\"\n$find2.Replacement.Text = \"programa\"\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$find2.Execute($find2.Text, $find2.MatchCase, $find2.MatchWholeWord, $find2.MatchWildcards, $null, $null, $find2.Forward, $find2.Wrap, $find2.Format, $find2.Replacement.Text, 2)\n"}
